# Apply the commit's changes to the "Sibirev I. V." worksheet.
# The edit fills in additional "5" scores across several rows and updates
# the frozen-pane view/selection to reflect where the instructor was
# working (around P8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 -----------------------------------------------------------
# New cell H4 (no prior style) gets a score of 5.
$ws.Range("H4").Value = 5

# --- Row 7 -------------------------------------------------------------
# D7 and F7 already carry style s="2"; just fill in the values.
$ws.Range("D7").Value = 5
$ws.Range("F7").Value = 5
# New cell H7 (no style) gets a score of 5.
$ws.Range("H7").Value = 5
# I7/J7/K7 already carry style s="8"; just fill in the values.
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 5

# --- Row 8 ---------------------------------------------------------
# F8 already carries style s="2"; just fill in the value.
$ws.Range("F8").Value = 5
# New cell H8 (no style) gets a score of 5.
$ws.Range("H8").Value = 5
# I8/J8/K8 already carry style s="8"; just fill in the values.
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 5
# New cell P8 (no style) gets a score of 5.
$ws.Range("P8").Value = 5

# --- Row 20 --------------------------------------------------------
# New cell H20 needs the same look as the other green "thick-left-border"
# cells (style index 19, e.g. V19) before getting its value.
$ws.Range("V19").Copy()
$ws.Range("H20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H20").Value = 5
$ws.Range("I20").Value = 5
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 5

# --- Row 25 ----------------------------------------------------------
# New cell G25 needs the "thick-both-borders" green look (style index 6,
# e.g. G30) before getting its value.
$ws.Range("G30").Copy()
$ws.Range("G25").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G25").Value = 5
# New cell H25 needs the same look as H20 (style index 19).
$ws.Range("V19").Copy()
$ws.Range("H25").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H25").Value = 5
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 5
$ws.Range("K25").Value = 5
# New cell L25 (no style) gets a score of 5.
$ws.Range("L25").Value = 5

# --- Row 26 ----------------------------------------------------------
$ws.Range("G30").Copy()
$ws.Range("G26").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G26").Value = 5
$ws.Range("V19").Copy()
$ws.Range("H26").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H26").Value = 5
$ws.Range("I26").Value = 5
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = 5
# New cell L26 (no style) gets a score of 5.
$ws.Range("L26").Value = 5

$excel.CutCopyMode = $false

# --- View state --------------------------------------------------------
# Move the active selection in the frozen bottom-right pane to P8, which
# is where this commit's edits were focused.
$ws.Range("P8").Select()

Write-Output "Applied AutoCommit edits"
